# Patton's Best - Events.xlsx
# "Added start of Evening Briefing"
#
# The e032 "No Combat" event description is revised: victory points are now
# only added to the After Action Report when the area is actually converted
# to US Control, instead of unconditionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

$newText = "<Bold>e032 No Combat</Bold> `n<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nIf converting territory to US Control, Victory points are added to the After Action Report `n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>."

$ws.Range("B33").Value = $newText
